$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("project_count") before the existing "reason" column,
# which shifts the old D column ("reason") to E.
$ws.Columns("D:D").Insert()

# --- Row 1: headers ---
$ws.Cells.Item(1, 4).Value = "project_count"
# E1 ("reason") already shifted into place by the column insert.

# --- Row 2 ---
$ws.Cells.Item(2, 3).Value = 73.92
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(2, 5).Value = "`nThe candidate has a score of 73.92 as they have experience in the relevant technologies required for the job, such as MongoDB, ReactJS, JavaScript, Web Development, NodeJS, as well as other related technologies such as Python, Django, Computer Vision, Image Processing. Their projects demonstrate their ability to work on developing web/mobile applications, feature development, scalability, and product enhancement."

# --- Row 3 ---
$ws.Cells.Item(3, 3).Value = 90.47
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = " `nReasoning: The candidate has demonstrated strong expertise in technologies relevant to the job description, such as MongoDB, ReactJS, JavaScript, Web Development, NodeJS, HTML, CSS, Socket.IO, WebRTC, Flutter, Dart, Firebase. Their projects also showcase a good understanding of the technologies and how to use them in combination to develop web/mobile applications. The candidate's score reflects their ability to develop high-quality applications in the aforementioned technologies."

# --- Row 4 ---
$ws.Cells.Item(4, 3).Value = 78.58
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = "`nReasoning: The candidate has demonstrated technical proficiency in MongoDB, ReactJS, JavaScript, Web Development, NodeJS, Django Rest Framework, Pytorch, Tensorflow, Keras, and Sklearn, which are the skills required for the job. The projects showcase the candidate's ability to design and develop high-quality web/mobile applications, integrate machine learning models with user interfaces, and develop automated pipelines for multi-model data analysis. All these qualities make the candidate a good fit for the job, justifying their score of 78.58."

# The multi-line "reason" text triggers an implicit row-height autofit with a
# pinned custom height; re-running AutoFit clears that pinned/custom height so
# the rows stay on the sheet's default (unset) height, matching the source file.
$ws.Rows("2:4").AutoFit()
